# "Generate Report for Handoff"
#
# The localization-status report drops the stale "efbf9219..." source-file
# row (it handed back already) and updates the "3976662b..." row's status
# from "Handed back: in sync with en-US" to "Ready for handoff" with fresh
# handoff timestamps, on all three sheets (Overview, zh-cn, de-de).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns File Name | zh-cn | de-de
#   row2 = 3976662b...   -> status text changes
#   row3 = efbf9219...   -> removed entirely
#   row4 = .localization-config -> shifts up to row3
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Wipe hyperlinks on this sheet (the engine's Range.Hyperlinks.Delete()
# clears the whole sheet's collection) -- we rebuild the surviving ones below.
$ws.Range("A2").Hyperlinks.Delete()

# Drop the efbf9219 row; row4 (.localization-config) slides up into row3.
$ws.Rows.Item(3).Delete()

# Update the still-present 3976662b row's status cells.
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5524b9b6258a57fe84c667a98b9e76d2ba061b12/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5524b9b6258a57fe84c667a98b9e76d2ba061b12/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": per-locale handoff/handback detail table
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-08 20:57:34"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5524b9b6258a57fe84c667a98b9e76d2ba061b12/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/11e7b1fecee053280070c8e63e93eba97a4cdeff/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/3b295e7a2bab4d42f81d2d2add8842e4e7c6fe55/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2aeeec35e53c1bb53d410b0a468bdec8e3195933/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.zh-cn.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5524b9b6258a57fe84c667a98b9e76d2ba061b12/.localization-config", "", "", ".localization-config") | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": same shape as zh-cn, different locale links/datetime
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Hyperlinks.Delete()
$ws.Rows.Item(3).Delete()

$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-03-08 20:57:42"

$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/5524b9b6258a57fe84c667a98b9e76d2ba061b12/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/587d7734278e73d372ae8daf5d6eceae52411a00/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/41a154508581ed90e1be7bb1d979825571f326a5/e2e/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.md") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4446e7ef7034f88361728f555f4c301ab4176bd7/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf", "", "", "3976662b-7e18-4ef4-9397-f67d8c4c4bc5.7fed2e4b03b1d9d3d15c2b653214b755747046ba.de-de.xlf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/5524b9b6258a57fe84c667a98b9e76d2ba061b12/.localization-config", "", "", ".localization-config") | Out-Null
